# LF Energy High Level Overview Deck - June 2025 update
#
# The canonical diff for this commit shows every <a:tbl> in the deck whose
# <a:tableStyleId> was {C1392B24-FFAB-4DD7-99E8-4D97128CA437} being
# re-pointed at {E6D9C5E6-9F9E-41A7-B0E6-88E029BFA50D}. There are exactly
# six such tables (one apiece on six slides). We walk every slide/shape,
# find any table using the old style id, and re-apply the new style via
# Table.ApplyStyle (Table.Style is read-only - PowerPoint raises an error
# if you try to assign it directly, and ApplyStyle is the supported verb).

$p = $ppt.ActivePresentation

$oldStyleId = "{C1392B24-FFAB-4DD7-99E8-4D97128CA437}"
$newStyleId = "{E6D9C5E6-9F9E-41A7-B0E6-88E029BFA50D}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)

        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style.Name -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
